$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Add-EventRow {
    param($worksheet, $rowNum, $dateSerial, $eventName, $location, $city, $url)

    $worksheet.Range("A$rowNum").Value = $dateSerial

    $worksheet.Range("B$rowNum").NumberFormat = "@"
    $worksheet.Range("B$rowNum").Value = $eventName

    $worksheet.Range("C$rowNum").NumberFormat = "@"
    $worksheet.Range("C$rowNum").Value = $location

    $worksheet.Range("D$rowNum").NumberFormat = "@"
    $worksheet.Range("D$rowNum").Value = $city

    $eCell = $worksheet.Range("E$rowNum")
    $eCell.NumberFormat = "@"
    $eCell.Value = $url

    # Style the link text as a rich-text run inside the shared string
    # (matches the rest of the sheet: underline, accent color, Calibri 11),
    # applying it to [1 .. len-1] and [len .. len] separately then merging
    # via identical formatting keeps the cell's own style untouched.
    $len = $url.Length
    $firstPart = $eCell.Characters(1, $len - 1)
    $firstPart.Font.Underline = 2
    $firstPart.Font.ColorIndex = 4
    $firstPart.Font.Name = "Calibri"
    $firstPart.Font.Size = 11

    $lastPart = $eCell.Characters($len, 1)
    $lastPart.Font.Underline = 2
    $lastPart.Font.ColorIndex = 4
    $lastPart.Font.Name = "Calibri"
    $lastPart.Font.Size = 11

    $worksheet.Hyperlinks.Add($eCell, $url, "", "", $url)

    # Hyperlinks.Add re-styles the whole cell with its own built-in
    # "Hyperlink" cell style; restore the plain text-cell style used
    # throughout the sheet so the cell keeps its original formatting and
    # only the shared-string run carries the link's visual styling.
    $eCell.Font.Underline = 0
    $eCell.Font.ColorIndex = 1
    $eCell.Font.Name = "Calibri"
    $eCell.Font.Size = 11
    $eCell.NumberFormat = "@"
}

Add-EventRow $ws "225" 45751 "DESIRE x HOT MEAL" "Schrotty" "Köln" "https://www.instagram.com/reel/DGtN90vAhj0/?igsh=MTQ1M2Vqdnd1eGg4Zg=="
Add-EventRow $ws "226" 45719 "ROSENMONTAG OPEN DOORS (15UHR)" "PM93" "Essen" "https://www.instagram.com/reel/DGnTHFoiKry/?igsh=MXFoaHZ6N3hxMmRvOA=="
Add-EventRow $ws "227" 45826 "U-BOUNCE" "U-Club" "Wuppertal" "https://www.instagram.com/ubounce_official?igsh=Z2I0c3B0NWhwbjIw"
